$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.966.68"
$ws.Range("E2").Value = "  +3.67%  "

$ws.Range("D3").Value = "3.265.31"
$ws.Range("E3").Value = "  +3.08%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'581.03"
$ws.Range("E5").Value = "  +1.94%  "

$ws.Range("D6").Value = "'182.24"
$ws.Range("E6").Value = "  +6.40%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("D9").Value = "3.265.23"
$ws.Range("E9").Value = "  +3.09%  "

$ws.Range("E10").Value = "  +7.36%  "

$ws.Range("D11").Value = "'6.73"
$ws.Range("E11").Value = "  +2.63%  "

$ws.Range("E12").Value = "  +6.51%  "

$ws.Range("D13").Value = "3.832.70"
$ws.Range("E13").Value = "  +3.14%  "

$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").Value = "'28.48"
$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("D16").Value = "67.946.45"
$ws.Range("E16").Value = "  +3.64%  "

$ws.Range("E17").Value = "  +3.88%  "

$ws.Range("D18").Value = "3.257.16"
$ws.Range("E18").Value = "  +2.69%  "

$ws.Range("E19").Value = "  +2.66%  "

$ws.Range("D20").Value = "'13.52"
$ws.Range("E20").Value = "  +5.09%  "

$ws.Range("D21").Value = "'375.78"
$ws.Range("E21").Value = "  +4.84%  "

$ws.Range("D22").Value = "'7.65"
$ws.Range("E22").Value = "  +5.40%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'71.18"
$ws.Range("E24").Value = "  +3.03%  "

$ws.Range("D25").Value = "'0.514"
$ws.Range("E25").Value = "  +4.07%  "

$ws.Range("E26").Value = "  +5.26%  "

$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("E28").Value = "  +2.20%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "'1.99"
$ws.Range("E30").Value = "  +3.51%  "

$ws.Range("D31").Value = "'5.68"
$ws.Range("E31").Value = "  +6.07%  "

$ws.Range("D32").Value = "'22.81"
$ws.Range("E32").Value = "  +4.11%  "

$ws.Range("E34").Value = "  +6.01%  "

$ws.Range("E35").Value = "  +4.86%  "

$ws.Range("E36").Value = "  +5.11%  "

$ws.Range("D37").Value = "'161.27"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").Value = "'0.850"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").Value = "'1.85"
$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("E40").Value = "  +11.05%  "

$ws.Range("D41").Value = "'26.79"
$ws.Range("E41").Value = "  +1.75%  "

$ws.Range("E42").Value = "  +11.22%  "

$ws.Range("E43").Value = "  +4.67%  "

$ws.Range("D44").Value = "'25.76"
$ws.Range("E44").Value = "  +7.49%  "

$ws.Range("D45").Value = "2.688.12"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("D46").Value = "'350.84"
$ws.Range("E46").Value = "  +7.00%  "

$ws.Range("D47").Value = "'40.85"
$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("E48").Value = "  +3.60%  "

$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("E50").Value = "  +5.87%  "

$ws.Range("E51").Value = "  +0.58%  "
